$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new numeric-looking text must remain literal text
# (mirrors how the source data was authored: column D stores free-form
# price strings, some of which look like plain decimals to Excel's parser).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '96.850.12'
$ws.Range("E2").Value = '  +0.68%  '

$ws.Range("D3").Value = '3.701.86'
$ws.Range("E3").Value = '  +3.84%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '242.71'
$ws.Range("E5").Value = '  +0.93%  '

$ws.Range("E6").Value = '  +19.89%  '

$ws.Range("D7").Value = '659.46'
$ws.Range("E7").Value = '  +0.97%  '

$ws.Range("D8").Value = '0.427'
$ws.Range("E8").Value = '  +5.75%  '

$ws.Range("E9").Value = '  +4.87%  '

$ws.Range("E10").Value = '  -0.04%  '

$ws.Range("D11").Value = '3.697.13'
$ws.Range("E11").Value = '  +3.86%  '

$ws.Range("D12").Value = '44.89'
$ws.Range("E12").Value = '  +4.05%  '

$ws.Range("E13").Value = '  +1.67%  '

$ws.Range("D14").Value = '6.54'
$ws.Range("E14").Value = '  +2.72%  '

$ws.Range("D15").Value = '4.395.65'
$ws.Range("E15").Value = '  +3.93%  '

$ws.Range("D16").Value = '96.823.49'
$ws.Range("E16").Value = '  +0.88%  '

$ws.Range("E17").Value = '  +1.94%  '

$ws.Range("D18").Value = '3.707.40'
$ws.Range("E18").Value = '  +4.02%  '

$ws.Range("D19").Value = '13.14'
$ws.Range("E19").Value = '  +5.20%  '

$ws.Range("E20").Value = '  +0.31%  '

$ws.Range("D21").Value = '18.50'
$ws.Range("E21").Value = '  +4.85%  '

$ws.Range("D22").Value = '0.549'
$ws.Range("E22").Value = '  +4.58%  '

$ws.Range("D23").Value = '513.96'
$ws.Range("E23").Value = '  +1.60%  '

$ws.Range("D24").Value = '3.45'
$ws.Range("E24").Value = '  +1.34%  '

$ws.Range("D25").Value = '0.0000211'
$ws.Range("E25").Value = '  +6.45%  '

$ws.Range("D26").Value = '6.90'
$ws.Range("E26").Value = '  +0.75%  '

$ws.Range("D27").Value = '101.61'
$ws.Range("E27").Value = '  +5.99%  '

$ws.Range("D28").Value = '13.07'
$ws.Range("E28").Value = '  +3.89%  '

$ws.Range("D29").Value = '0.169'
$ws.Range("E29").Value = '  +12.72%  '

$ws.Range("D30").Value = '3.04'
$ws.Range("E30").Value = '  +2.59%  '

$ws.Range("D31").Value = '12.00'
$ws.Range("E31").Value = '  +5.72%  '

$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.12%  '

$ws.Range("E33").Value = '  +2.18%  '

$ws.Range("D34").Value = '33.30'
$ws.Range("E34").Value = '  +6.36%  '

$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.25%  '

$ws.Range("D36").Value = '0.593'
$ws.Range("E36").Value = '  +5.71%  '

$ws.Range("E37").Value = '  +6.52%  '

$ws.Range("D38").Value = '612.71'
$ws.Range("E38").Value = '  -0.43%  '

$ws.Range("D39").Value = '8.70'
$ws.Range("E39").Value = '  -0.29%  '

$ws.Range("D40").Value = '42.59'
$ws.Range("E40").Value = '  +27.49%  '

$ws.Range("E41").Value = '  +6.23%  '

$ws.Range("D42").Value = '0.966'
$ws.Range("E42").Value = '  +7.50%  '

$ws.Range("E43").Value = '  +7.26%  '

$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("D45").Value = '6.11'
$ws.Range("E45").Value = '  +7.65%  '

$ws.Range("D46").Value = '0.0442'
$ws.Range("E46").Value = '  +5.36%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.418'
$ws.Range("E47").Value = '  +22.66%  '

$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '2.30'
$ws.Range("E48").Value = '  +1.37%  '

$ws.Range("E49").Value = '  +0.35%  '

$ws.Range("D50").Value = '8.57'
$ws.Range("E50").Value = '  +5.21%  '

$ws.Range("D51").Value = '54.51'
$ws.Range("E51").Value = '  +3.25%  '
